# Fixed update to excel issue
# - Renames the "Requested quantity" headers on the existing sheets so they
#   carry distinguishing names per source table.
# - Adds a new "PO Forecast" worksheet (appended after the existing sheets)
#   containing the Prophet-style forecast columns: ds, PO_Forecast,
#   yhat_lower, yhat_upper.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" header on the weekly sheet ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Rename the "Requested quantity" header on the monthly sheet ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the outline + page-margin conventions used by the other sheets.
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/centered header style already used on the other sheets.
$wsWeekly.Range("A1:B1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null

# Forecast data rows (ds serial, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45109.99999999999, 10, -10.09419099274506, 31.06959981083741),
    @(45116.99999999999, 10, -10.45922346365151, 30.00429069574198),
    @(45158.99999999999, 11, -9.220621540417106, 31.46872826292023),
    @(45200.99999999999, 12, -8.425437367480338, 31.61269463765984),
    @(45221.99999999999, 13, -8.192366793242243, 33.21015405408936),
    @(45487.99999999999, 18, -1.002623487348484, 37.90545454395797),
    @(45494.99999999999, 18, -1.951780724505699, 37.70794771835185),
    @(45501.99999999999, 18, -1.927817898801446, 38.49715633655051),
    @(45508.99999999999, 18, -2.345806566267158, 38.93752618106274),
    @(45515.99999999999, 18, -2.809502225123798, 39.54497006951394),
    @(45522.99999999999, 18, -2.740885723889634, 38.07945509468755),
    @(45529.99999999999, 19, -1.66344440181685, 40.02040989129679),
    @(45536.99999999999, 19, -1.32379866790631, 39.41605477539296),
    @(45543.99999999999, 19, -1.491663790209314, 39.60754192357727),
    @(45550.99999999999, 19, -0.2893790244826248, 39.64022998057851),
    @(45557.99999999999, 19, -0.8279046360947149, 39.25840072921974),
    @(45564.99999999999, 19, -1.69725816741057, 39.06037295814338),
    @(45585.99999999999, 20, -0.5829596408297242, 42.08125621834862),
    @(45592.99999999999, 20, -0.7127842026146555, 39.96904793143984),
    @(45599.99999999999, 20, 1.287081773957452, 40.33250682075049),
    @(45606.99999999999, 20, 0.5152123251515992, 38.38134303414625),
    @(45613.99999999999, 20, -1.049012903050061, 39.77193087731921),
    @(45620.99999999999, 20, -0.3672159860896225, 41.86416697295713),
    @(45627.99999999999, 20, -0.3849199068045632, 39.99491414573833),
    @(45634.99999999999, 21, 0.4066433206818195, 40.57222139501176),
    @(45641.99999999999, 21, 0.4790644642033959, 41.21764404240746)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the date/time number-format style used for column A ("ds") on the
# other sheets so the forecast dates render the same way.
$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A27").PasteSpecial(-4122) | Out-Null

$wsForecast.Range("A1").Select() | Out-Null
